{"js": "// Update the date line and the twenty-five two-digit multiplication\n// answers in the table to the new values from the commit.\nconst replacements = [\n  [\"2023-08-16 Wednesday\", \"2023-08-17 Thursday\"],\n  [\"95\u00d723=2185\", \"95\u00d787=8265\"],\n  [\"43\u00d768=2924\", \"40\u00d727=1080\"],\n  [\"29\u00d722=638\", \"77\u00d748=3696\"],\n  [\"60\u00d766=3960\", \"65\u00d740=2600\"],\n  [\"48\u00d733=1584\", \"31\u00d763=1953\"],\n  [\"19\u00d743=817\", \"98\u00d749=4802\"],\n  [\"81\u00d799=8019\", \"28\u00d760=1680\"],\n  [\"85\u00d752=4420\", \"14\u00d739=546\"],\n  [\"72\u00d768=4896\", \"19\u00d747=893\"],\n  [\"70\u00d775=5250\", \"25\u00d798=2450\"],\n  [\"80\u00d797=7760\", \"51\u00d758=2958\"],\n  [\"37\u00d713=481\", \"36\u00d764=2304\"],\n  [\"13\u00d770=910\", \"25\u00d731=775\"],\n  [\"32\u00d777=2464\", \"88\u00d717=1496\"],\n  [\"17\u00d792=1564\", \"40\u00d728=1120\"],\n  [\"29\u00d732=928\", \"75\u00d737=2775\"],\n  [\"43\u00d741=1763\", \"23\u00d749=1127\"],\n  [\"40\u00d741=1640\", \"97\u00d742=4074\"],\n  [\"75\u00d754=4050\", \"82\u00d773=5986\"],\n  [\"38\u00d750=1900\", \"86\u00d744=3784\"],\n  [\"60\u00d757=3420\", \"82\u00d783=6806\"],\n  [\"89\u00d760=5340\", \"99\u00d733=3267\"],\n  [\"95\u00d715=1425\", \"34\u00d722=748\"],\n  [\"26\u00d734=884\", \"80\u00d755=4400\"],\n  [\"67\u00d764=4288\", \"86\u00d721=1806\"],\n];\n\nfor (const [oldText, newText] of replacements) {\n  const results = context.document.body.search(oldText, { matchCase: true });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, \"Replace\");\n  }\n  await context.sync();\n}\n", "ps1": "# Update the date line and the twenty-five two-digit multiplication\n# answers in the table to the new values from the commit.\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"2023-08-16 Wednesday\", \"2023-08-17 Thursday\"),\n    @(\"95\u00d723=2185\", \"95\u00d787=8265\"),\n    @(\"43\u00d768=2924\", \"40\u00d727=1080\"),\n    @(\"29\u00d722=638\", \"77\u00d748=3696\"),\n    @(\"60\u00d766=3960\", \"65\u00d740=2600\"),\n    @(\"48\u00d733=1584\", \"31\u00d763=1953\"),\n    @(\"19\u00d743=817\", \"98\u00d749=4802\"),\n    @(\"81\u00d799=8019\", \"28\u00d760=1680\"),\n    @(\"85\u00d752=4420\", \"14\u00d739=546\"),\n    @(\"72\u00d768=4896\", \"19\u00d747=893\"),\n    @(\"70\u00d775=5250\", \"25\u00d798=2450\"),\n    @(\"80\u00d797=7760\", \"51\u00d758=2958\"),\n    @(\"37\u00d713=481\", \"36\u00d764=2304\"),\n    @(\"13\u00d770=910\", \"25\u00d731=775\"),\n    @(\"32\u00d777=2464\", \"88\u00d717=1496\"),\n    @(\"17\u00d792=1564\", \"40\u00d728=1120\"),\n    @(\"29\u00d732=928\", \"75\u00d737=2775\"),\n    @(\"43\u00d741=1763\", \"23\u00d749=1127\"),\n    @(\"40\u00d741=1640\", \"97\u00d742=4074\"),\n    @(\"75\u00d754=4050\", \"82\u00d773=5986\"),\n    @(\"38\u00d750=1900\", \"86\u00d744=3784\"),\n    @(\"60\u00d757=3420\", \"82\u00d783=6806\"),\n    @(\"89\u00d760=5340\", \"99\u00d733=3267\"),\n    @(\"95\u00d715=1425\", \"34\u00d722=748\"),\n    @(\"26\u00d734=884\", \"80\u00d755=4400\"),\n    @(\"67\u00d764=4288\", \"86\u00d721=1806\")\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Replacement.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null\n}\n"}
